# Apply selection updates and the 2024-sheet "Annual Return" value
# refresh across the workbook. Active sheet stays "2024" (5th, 0-indexed)
# to match the committed workbook state.

$wb = $excel.ActiveWorkbook

# --- 2019 ---
$ws1 = $wb.Worksheets.Item("2019")
$ws1.Activate()
$ws1.Range("C41").Select()

# --- 2020 ---
$ws2 = $wb.Worksheets.Item("2020")
$ws2.Activate()
$ws2.Range("D55").Select()

# --- 2021 ---
$ws3 = $wb.Worksheets.Item("2021")
$ws3.Activate()
$ws3.Range("C58").Select()

# --- 2022 ---
$ws4 = $wb.Worksheets.Item("2022")
$ws4.Activate()
$ws4.Range("C15").Select()

# --- 2023 ---
$ws5 = $wb.Worksheets.Item("2023")
$ws5.Activate()
$ws5.Range("B29").Select()

# --- 2024 (updated Annual Return values in column B) ---
$ws6 = $wb.Worksheets.Item("2024")
$ws6.Activate()

$ws6.Range("B2").Value = -0.11145510835913321
$ws6.Range("B3").Value = 0.1080213849698837
$ws6.Range("B4").Value = 0.1064854759850451
$ws6.Range("B5").Value = -0.054995970991135963
$ws6.Range("B6").Value = 0.065626681011295762
$ws6.Range("B7").Value = 0.098613251155623027
$ws6.Range("B8").Value = 0.16674197921373629
$ws6.Range("B9").Value = 0.06572144126676438
$ws6.Range("B10").Value = 0.1024925198965982
$ws6.Range("B11").Value = 0.15737473535638569
$ws6.Range("B12").Value = 0.20844055584148211
$ws6.Range("B13").Value = -0.1035707779194754
$ws6.Range("B14").Value = 0.1512388162422578
$ws6.Range("B15").Value = -0.01020408163265518
$ws6.Range("B16").Value = -0.099429115128448764
$ws6.Range("B17").Value = 0.33494872563012312
$ws6.Range("B18").Value = -0.04953497775980531
$ws6.Range("B19").Value = 0.062525375558264429
$ws6.Range("B20").Value = 0.13551401869158911
$ws6.Range("B21").Value = 0.053673627223510767
$ws6.Range("B22").Value = -0.077795104835467588
$ws6.Range("B23").Value = 0.15661252900232039
$ws6.Range("B24").Value = 0.1960000000000004
$ws6.Range("B25").Value = 0.1978319783197833
$ws6.Range("B26").Value = 0.15695346795434631
$ws6.Range("B27").Value = 0.02507712210170188
$ws6.Range("B28").Value = 0.15339902204943151
$ws6.Range("B29").Value = 0.2442384769539061
$ws6.Range("B30").Value = 0.23721954948069121
$ws6.Range("B31").Value = 0.081820050709243919
$ws6.Range("B32").Value = 0.1209915539380626
$ws6.Range("B33").Value = 0.058163720101025218
$ws6.Range("B34").Value = 0.088006986221619732
$ws6.Range("B35").Value = -0.077976817702844148
$ws6.Range("B36").Value = 0.11336982017200881
$ws6.Range("B37").Value = -0.054995970991135963
$ws6.Range("B38").Value = -0.1017942145734156
$ws6.Range("B39").Value = 0.1949567181031244
$ws6.Range("B40").Value = -0.08941485864562837
$ws6.Range("B41").Value = 0.042255511588468853
$ws6.Range("B42").Value = 0.2035928143712564
$ws6.Range("B43").Value = 0.12836624775583449
$ws6.Range("B44").Value = 0.2345554195711155
$ws6.Range("B45").Value = 0.072084160807257769
$ws6.Range("B46").Value = 0.1041515517936311
$ws6.Range("B47").Value = 0.23921683734878291
$ws6.Range("B48").Value = -0.11759504862953141
$ws6.Range("B49").Value = 0.23151645979492949
$ws6.Range("B50").Value = 0.062525375558264429
$ws6.Range("B51").Value = 0.13508260447036119
$ws6.Range("B52").Value = 0.044333149601808897
$ws6.Range("B53").Value = 0.25660226561956878
$ws6.Range("B54").Value = 0.2896855398598932
$ws6.Range("B55").Value = 0.2605398675796502
$ws6.Range("B56").Value = 0.1613361762615482
$ws6.Range("B57").Value = 0.11354817140878801
$ws6.Range("B58").Value = 0.1122944452457038
$ws6.Range("B60").Value = 0.28316197539187637
$ws6.Range("B61").Value = 0.10280569514237679
$ws6.Range("B62").Value = 0.2099832211711701
$ws6.Range("B63").Value = 0.26649041375039179
$ws6.Range("B64").Value = 0.43636363636363562
$ws6.Range("B65").Value = 0.19541875447387369
$ws6.Range("B66").Value = 0.33487677537260557
$ws6.Range("B67").Value = 0.2512421815631023
$ws6.Range("B68").Value = 0.2209185561630633
$ws6.Range("B69").Value = 0.27036245229021172
$ws6.Range("B70").Value = 0.24389314105734991
$ws6.Range("B71").Value = 0.21189206832771279
$ws6.Range("B72").Value = 0.22918654464200361
$ws6.Range("B73").Value = 0.23914592902533191
$ws6.Range("B74").Value = 0.14163017671182329
$ws6.Range("B75").Value = 0.1236887844102068
$ws6.Range("B76").Value = 0.1598281700848361
$ws6.Range("B77").Value = 0.18489055269588769
$ws6.Range("B78").Value = 0.20438930024681559
$ws6.Range("B79").Value = 0.1211136973086855

$ws6.Range("P19").Select()
